# Update "想去人数" (wanted-to-go count) values in the "展览" sheet
# and the corresponding rows in the "全部类型" sheet, matching the
# regenerated data snapshot described in the commit.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAll = $wb.Worksheets.Item("全部类型")

# Row => new F-value updates for "展览" sheet
$exhibitionUpdates = @{
    5  = 73
    9  = 8731
    10 = 807
    12 = 1146
    13 = 984
    14 = 110
    21 = 1029
}

foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Row => new F-value updates for "全部类型" sheet (same events, offset by one row)
$allTypesUpdates = @{
    6  = 73
    11 = 8731
    12 = 807
    14 = 1146
    15 = 984
    16 = 110
    23 = 1029
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allTypesUpdates[$row]
}

$wb.Save()
